$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.122.12"
$ws.Range("E2").Value = "  +3.51%  "

$ws.Range("D3").Value = "1.577.10"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("E4").Value = "  -1.12%  "

$ws.Range("D5").Value = "'212.91"
$ws.Range("E5").Value = "  +0.62%  "

$ws.Range("D6").Value = "'0.493"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("E7").Value = "  -1.07%  "

$ws.Range("D8").Value = "'23.27"
$ws.Range("E8").Value = "  +5.73%  "

$ws.Range("E9").Value = "  +0.64%  "

$ws.Range("D10").Value = "'0.0599"
$ws.Range("E10").Value = "  -0.17%  "

$ws.Range("D11").Value = "'0.0882"
$ws.Range("E11").Value = "  +1.74%  "

$ws.Range("D12").Value = "1.801.95"
$ws.Range("E12").Value = "  +0.49%  "

$ws.Range("D13").Value = "1.575.76"
$ws.Range("E13").Value = "  +0.44%  "

$ws.Range("E14").Value = "  -0.68%  "

$ws.Range("D15").Value = "'0.524"
$ws.Range("E15").Value = "  +0.84%  "

$ws.Range("D16").Value = "28.079.42"
$ws.Range("E16").Value = "  +3.40%  "

$ws.Range("D17").Value = "'63.65"
$ws.Range("E17").Value = "  +2.10%  "

$ws.Range("D18").Value = "'229.12"
$ws.Range("E18").Value = "  +6.11%  "

$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("D20").Value = "'7.46"
$ws.Range("E20").Value = "  +0.69%  "

$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.98%  "

$ws.Range("E22").Value = "  -0.71%  "

$ws.Range("D23").Value = "'9.33"
$ws.Range("E23").Value = "  +1.18%  "

$ws.Range("E24").Value = "  +0.29%  "

$ws.Range("D25").Value = "'152.17"
$ws.Range("E25").Value = "  -1.16%  "

$ws.Range("D26").Value = "'15.24"
$ws.Range("E26").Value = "  +1.05%  "

$ws.Range("D27").Value = "'6.59"

$ws.Range("E28").Value = "  +0.28%  "

$ws.Range("E29").Value = "  -0.99%  "

$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").Value = "'0.0474"
$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("E32").Value = "  -0.55%  "

$ws.Range("E33").Value = "  -1.75%  "

$ws.Range("D34").Value = "1.415.60"
$ws.Range("E34").Value = "  -2.37%  "

$ws.Range("E35").Value = "  -1.31%  "

$ws.Range("E36").Value = "  -4.89%  "

$ws.Range("E37").Value = "  -1.32%  "

$ws.Range("E38").Value = "  -0.21%  "

$ws.Range("D39").Value = "'0.541"
$ws.Range("E39").Value = "  +1.14%  "

$ws.Range("E40").Value = "  +5.83%  "

$ws.Range("D41").Value = "'0.807"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "  -1.12%  "

$ws.Range("D43").Value = "'5.68"
$ws.Range("E43").Value = "  -2.40%  "

$ws.Range("D44").Value = "'0.972"
$ws.Range("E44").Value = "  -2.78%  "

$ws.Range("E45").Value = "  +4.98%  "

$ws.Range("D46").Value = "'63.86"
$ws.Range("E46").Value = "  -1.25%  "

$ws.Range("D47").Value = "1.715.12"
$ws.Range("E47").Value = "  +0.61%  "

$ws.Range("D48").Value = "'87.00"
$ws.Range("E48").Value = "  +1.16%  "

$ws.Range("E49").Value = "  +2.13%  "

$ws.Range("D50").Value = "'0.0525"
$ws.Range("E50").Value = "  +0.88%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0943"
$ws.Range("E51").Value = "  -1.66%  "

